$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '67.777.87'
$ws.Range('E2').Value = '  -0.09%  '
$ws.Range('D3').Value = '3.236.65'
$ws.Range('E3').Value = '  -0.48%  '
$ws.Range('E4').Value = '  +0.01%  '
$ws.Range('D5').Value = '579.92'
$ws.Range('E5').Value = '  -0.63%  '
$ws.Range('D6').Value = '183.65'
$ws.Range('E6').Value = '  +0.22%  '
$ws.Range('E7').Value = '  +0.01%  '
$ws.Range('E8').Value = '  +0.26%  '
$ws.Range('E9').Value = '  -4.00%  '
$ws.Range('D10').Value = '6.59'
$ws.Range('E10').Value = '  -1.32%  '
$ws.Range('E11').Value = '  -0.32%  '
$ws.Range('D12').Value = '3.799.00'
$ws.Range('E12').Value = '  -0.41%  '
$ws.Range('E13').Value = '  +0.04%  '
$ws.Range('D14').Value = '27.57'
$ws.Range('E14').Value = '  -3.48%  '
$ws.Range('D15').Value = '67.814.31'
$ws.Range('E15').Value = '  -0.04%  '
$ws.Range('E16').Value = '  -1.89%  '
$ws.Range('D17').Value = '3.243.26'
$ws.Range('E17').Value = '  -0.33%  '
$ws.Range('E18').Value = '  -1.17%  '
$ws.Range('D19').Value = '13.43'
$ws.Range('E19').Value = '  -1.10%  '
$ws.Range('D20').Value = '395.93'
$ws.Range('E20').Value = '  +3.86%  '
$ws.Range('E21').Value = '  -1.33%  '
$ws.Range('D22').Value = '0.999'
$ws.Range('E22').Value = '  -0.27%  '
$ws.Range('D23').Value = '71.00'
$ws.Range('E23').Value = '  -0.31%  '
$ws.Range('E24').Value = '  +0.03%  '
$ws.Range('E25').Value = '  -1.75%  '
$ws.Range('E26').Value = '  +2.69%  '
$ws.Range('D27').Value = '9.52'
$ws.Range('E27').Value = '  -3.17%  '
$ws.Range('E28').Value = '  +0.09%  '
$ws.Range('E29').Value = '  -1.46%  '
$ws.Range('E30').Value = '  -1.91%  '
$ws.Range('D31').Value = '22.65'
$ws.Range('E31').Value = '  -1.07%  '
$ws.Range('E32').Value = '  -2.37%  '
$ws.Range('E33').Value = '  -0.95%  '
$ws.Range('D35').Value = '161.69'
$ws.Range('E35').Value = '  +0.27%  '
$ws.Range('E36').Value = '  -4.44%  '
$ws.Range('E37').Value = '  +1.90%  '
$ws.Range('B38').Value = 'Mantle'
$ws.Range('C38').Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range('D38').Value = '0.808'
$ws.Range('E38').Value = '  -3.42%  '
$ws.Range('B39').Value = 'EnergySwap'
$ws.Range('C39').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D39').Value = '26.36'
$ws.Range('E39').Value = '  -0.62%  '
$ws.Range('E40').Value = '  -1.10%  '
$ws.Range('D41').Value = '6.45'
$ws.Range('E41').Value = '  -3.81%  '
$ws.Range('D42').Value = '41.10'
$ws.Range('E42').Value = '  -0.46%  '
$ws.Range('D43').Value = '2.45'
$ws.Range('E43').Value = '  -4.91%  '
$ws.Range('D44').Value = '0.0682'
$ws.Range('E44').Value = '  -0.75%  '
$ws.Range('B45').Value = 'Maker'
$ws.Range('C45').Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range('D45').Value = '2.603.36'
$ws.Range('E45').Value = '  -1.01%  '
$ws.Range('B46').Value = 'InjectiveProtocol'
$ws.Range('C46').Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range('D46').Value = '24.90'
$ws.Range('E46').Value = '  -1.97%  '
$ws.Range('D47').Value = '334.99'
$ws.Range('E47').Value = '  -3.48%  '
$ws.Range('D48').Value = '0.0277'
$ws.Range('E48').Value = '  -2.26%  '
$ws.Range('D49').Value = '6.26'
$ws.Range('E49').Value = '  +0.80%  '
$ws.Range('E50').Value = '  -2.06%  '
$ws.Range('D51').Value = '30.90'
$ws.Range('E51').Value = '  +0.26%  '
